$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates: volume/number and week-of date range ---
$ws.Range("A8").Value = "Volume 32   Number  46"
$ws.Range("C9").Value = "Report Covering the Week  11/10/2025  Through  11/16/2025"

# --- Weekly crime-complaint table updates (rows 14-30) ---
# Row 14
$ws.Range("N14").Value2 = -87.301587301587
# Row 15
$ws.Range("D15").Value2 = 2
$ws.Range("G15").Value2 = 7
$ws.Range("H15").Value2 = -42.857142857142
$ws.Range("J15").Value2 = 40
$ws.Range("K15").Value2 = -15
$ws.Range("N15").Value2 = -58.536585365853
# Row 16
$ws.Range("C16").Value2 = 6
$ws.Range("D16").Value2 = 16
$ws.Range("E16").Value2 = -62.5
$ws.Range("F16").Value2 = 35
$ws.Range("G16").Value2 = 59
$ws.Range("H16").Value2 = -40.677966101694
$ws.Range("I16").Value2 = 410
$ws.Range("J16").Value2 = 511
$ws.Range("K16").Value2 = -19.765166340508
$ws.Range("L16").Value2 = -14.760914760914
$ws.Range("M16").Value2 = -5.092592592592
$ws.Range("N16").Value2 = -77.120535714285
# Row 17
$ws.Range("C17").Value2 = 13
$ws.Range("D17").Value2 = 16
$ws.Range("E17").Value2 = -18.75
$ws.Range("F17").Value2 = 65
$ws.Range("G17").Value2 = 84
$ws.Range("H17").Value2 = -22.619047619047
$ws.Range("I17").Value2 = 910
$ws.Range("J17").Value2 = 939
$ws.Range("K17").Value2 = -3.088391906283
$ws.Range("L17").Value2 = -6.088751289989
$ws.Range("M17").Value2 = 86.094069529652
$ws.Range("N17").Value2 = -25.287356321839
# Row 18
$ws.Range("C18").Value2 = 5
$ws.Range("E18").Value2 = -28.571428571428
$ws.Range("G18").Value2 = 28
$ws.Range("H18").Value2 = -39.285714285714
$ws.Range("I18").Value2 = 208
$ws.Range("J18").Value2 = 264
$ws.Range("K18").Value2 = -21.212121212121
$ws.Range("L18").Value2 = -48
$ws.Range("M18").Value2 = -10.729613733905
$ws.Range("N18").Value2 = -88.732394366197
# Row 19
$ws.Range("C19").Value2 = 13
$ws.Range("D19").Value2 = 20
$ws.Range("E19").Value2 = -35
$ws.Range("F19").Value2 = 59
$ws.Range("G19").Value2 = 89
$ws.Range("H19").Value2 = -33.707865168539
$ws.Range("I19").Value2 = 665
$ws.Range("J19").Value2 = 803
$ws.Range("K19").Value2 = -17.185554171855
$ws.Range("L19").Value2 = -6.862745098039
$ws.Range("M19").Value2 = 98.507462686567
$ws.Range("N19").Value2 = -7.766990291262
# Row 20
$ws.Range("C20").Value2 = 7
$ws.Range("D20").Value2 = 3
$ws.Range("E20").Value2 = 133.333333333333
$ws.Range("F20").Value2 = 18
$ws.Range("G20").Value2 = 24
$ws.Range("H20").Value2 = -25
$ws.Range("I20").Value2 = 289
$ws.Range("J20").Value2 = 318
$ws.Range("K20").Value2 = -9.119496855345
$ws.Range("L20").Value2 = -31.190476190476
$ws.Range("M20").Value2 = 57.065217391304
$ws.Range("N20").Value2 = -75.086206896551
# Row 21
$ws.Range("C21").Value2 = 44
$ws.Range("D21").Value2 = 64
$ws.Range("E21").Value2 = -31.25
$ws.Range("G21").Value2 = 291
$ws.Range("H21").Value2 = -31.958762886597
$ws.Range("I21").Value2 = 2524
$ws.Range("J21").Value2 = 2888
$ws.Range("K21").Value2 = -12.603878116343
$ws.Range("L21").Value2 = -16.809492419248
$ws.Range("M21").Value2 = 47.257876312718
$ws.Range("N21").Value2 = -63.324614937518
# Row 22
$ws.Range("C22").Value2 = 1
$ws.Range("C22").NumberFormat = "#,##0"
$ws.Range("D22").Value2 = 1
$ws.Range("D22").NumberFormat = "#,##0"
$ws.Range("E22").Value2 = 0
$ws.Range("E22").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F22").Value2 = 2
$ws.Range("G22").Value2 = 3
$ws.Range("H22").Value2 = -33.333333333333
$ws.Range("I22").Value2 = 37
$ws.Range("J22").Value2 = 45
$ws.Range("K22").Value2 = -17.777777777777
$ws.Range("L22").Value2 = -2.631578947368
$ws.Range("M22").Value2 = 0
# Row 23
$ws.Range("D23").Value2 = 2
$ws.Range("D23").NumberFormat = "#,##0"
$ws.Range("E23").Value2 = -100
$ws.Range("E23").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F23").Value2 = 3
$ws.Range("G23").Value2 = 5
$ws.Range("H23").Value2 = -40
$ws.Range("J23").Value2 = 68
$ws.Range("K23").Value2 = -44.117647058823
$ws.Range("L23").Value2 = -45.714285714285
# Row 24
$ws.Range("C24").Value2 = 39
$ws.Range("D24").Value2 = 36
$ws.Range("E24").Value2 = 8.333333333333
$ws.Range("F24").Value2 = 106
$ws.Range("G24").Value2 = 145
$ws.Range("H24").Value2 = -26.896551724137
$ws.Range("I24").Value2 = 1585
$ws.Range("J24").Value2 = 1708
$ws.Range("K24").Value2 = -7.201405152224
$ws.Range("L24").Value2 = -14.55525606469
$ws.Range("M24").Value2 = 31.100082712985
# Row 25
$ws.Range("C25").Value2 = 11
$ws.Range("D25").Value2 = 12
$ws.Range("E25").Value2 = -8.333333333333
$ws.Range("F25").Value2 = 31
$ws.Range("G25").Value2 = 65
$ws.Range("H25").Value2 = -52.307692307692
$ws.Range("I25").Value2 = 586
$ws.Range("J25").Value2 = 776
$ws.Range("K25").Value2 = -24.484536082474
$ws.Range("L25").Value2 = -39.211618257261
# Row 26
$ws.Range("C26").Value2 = 21
$ws.Range("D26").Value2 = 21
$ws.Range("E26").Value2 = 0
$ws.Range("F26").Value2 = 93
$ws.Range("G26").Value2 = 107
$ws.Range("H26").Value2 = -13.084112149532
$ws.Range("I26").Value2 = 1142
$ws.Range("J26").Value2 = 1269
$ws.Range("K26").Value2 = -10.007880220646
$ws.Range("L26").Value2 = 4.387568555758
$ws.Range("M26").Value2 = 3.068592057761
# Row 27
$ws.Range("C27").Value2 = 1
$ws.Range("C27").NumberFormat = "#,##0"
$ws.Range("D27").Value2 = 2
$ws.Range("E27").Value2 = -50
$ws.Range("G27").Value2 = 8
$ws.Range("H27").Value2 = -25
$ws.Range("I27").Value2 = 53
$ws.Range("J27").Value2 = 60
$ws.Range("K27").Value2 = -11.666666666666
$ws.Range("L27").Value2 = -13.114754098360
# Row 28
$ws.Range("D28").Value2 = 3
$ws.Range("E28").Value2 = 100
$ws.Range("F28").Value2 = 16
$ws.Range("G28").Value2 = 17
$ws.Range("H28").Value2 = -5.882352941176
$ws.Range("I28").Value2 = 122
$ws.Range("J28").Value2 = 123
$ws.Range("K28").Value2 = -0.813008130081
$ws.Range("L28").Value2 = 2.521008403361
# Row 29
$ws.Range("C29").Value2 = 1
$ws.Range("C29").NumberFormat = "#,##0"
$ws.Range("D29").Value2 = 1
$ws.Range("D29").NumberFormat = "#,##0"
$ws.Range("E29").Value2 = 0
$ws.Range("E29").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F29").Value2 = 2
$ws.Range("H29").Value2 = 0
$ws.Range("I29").Value2 = 32
$ws.Range("J29").Value2 = 57
$ws.Range("K29").Value2 = -43.859649122807
$ws.Range("L29").Value2 = -36
$ws.Range("M29").Value2 = -13.513513513513
$ws.Range("N29").Value2 = -80.487804878048
# Row 30
$ws.Range("C30").Value2 = 1
$ws.Range("C30").NumberFormat = "#,##0"
$ws.Range("D30").Value2 = 1
$ws.Range("D30").NumberFormat = "#,##0"
$ws.Range("E30").Value2 = 0
$ws.Range("E30").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F30").Value2 = 2
$ws.Range("H30").Value2 = 0
$ws.Range("I30").Value2 = 25
$ws.Range("J30").Value2 = 42
$ws.Range("K30").Value2 = -40.476190476190
$ws.Range("L30").Value2 = -37.5
$ws.Range("M30").Value2 = -16.666666666666
$ws.Range("N30").Value2 = -83.108108108108
